{"js": "// hotfix: the date for the januari invoice\n// The invoice's \"Date of Invoice\" value changes from 12-01-2026 to 20-01-2026.\nconst body = context.document.body;\n\nconst results = body.search(\"12-01-2026\", { matchCase: true, matchWholeWord: false });\nresults.load(\"items/text\");\nawait context.sync();\n\nif (results.items.length === 0) {\n  throw new Error(\"Could not find the invoice date '12-01-2026' to update.\");\n}\n\n// Replace every match (expected: exactly one, the \"Date of Invoice\" cell).\nfor (const match of results.items) {\n  match.insertText(\"20-01-2026\", \"Replace\");\n}\n\nawait context.sync();\n", "ps1": "# hotfix: the date for the januari invoice\n# The invoice's \"Date of Invoice\" value changes from 12-01-2026 to 20-01-2026.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"12-01-2026\"\n$find.Replacement.Text = \"20-01-2026\"\n$find.Forward = $true\n$find.Wrap = 1            # wdFindContinue\n$find.Format = $false\n$find.MatchCase = $true\n$find.MatchWholeWord = $false\n$find.MatchWildcards = $false\n\n$found = $find.Execute(\n  $find.Text,\n  $false,\n  $true,\n  $false,\n  $false,\n  $false,\n  $true,\n  1,\n  $false,\n  $find.Replacement.Text,\n  2\n)\n\nif (-not $found) {\n  throw \"Could not find the invoice date '12-01-2026' to update.\"\n}\n"}
